$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124 - this shifts every existing row
# (124..153) down by one, matching the weekly-data-prepend pattern in
# the diff (old row 124 data reappears as new row 125, ..., old row 153
# data reappears as new row 154).
$ws.Rows.Item(124).Insert()

# Populate the freshly inserted row 124 with the new weekly entry.
$ws.Range("A124").Value = 10
$ws.Range("B124").Value = "Vega Modelo de Temuco"
$ws.Range("C124").Value = "La Araucanía"
$ws.Range("D124").Value = 44511
$ws.Range("E124").Value = 9
$ws.Range("F124").Value = 100112005
$ws.Range("G124").Value = "Puerro"
$ws.Range("H124").Value = "Azul de Maquehue"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 70
$ws.Range("K124").Value = 8000
$ws.Range("L124").Value = 8000
$ws.Range("M124").Value = 8000
$ws.Range("N124").Value = '$/docena de paquetes'
$ws.Range("O124").Value = "Provincia de Cautín"
$ws.Range("P124").Value = 667
$ws.Range("Q124").Value = 12
$ws.Range("R124").Value = "Hortaliza"
